$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A20").Value = "Appellant statement received"
$ws.Range("A21").Value = "Appellant statement reference number"

$ws.Range("A21").Select()
